$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: insert a new paragraph (reproducibility note) immediately before
# the "Advanced Audit Policy..." paragraph (i.e. right after the
# "Experimental environment" Heading2 paragraph).
#
# We insert the sentence as plain text at the very start of the "Advanced
# Audit Policy" paragraph (its first run carries no rPr, so the inserted
# text picks up no formatting either) and then split it off into its own
# paragraph with InsertParagraphAfter - this keeps both the new paragraph
# and the following one free of any pPr/rPr noise, matching a plain
# "<w:p><w:r><w:t>...</w:t></w:r></w:p>".
# ---------------------------------------------------------------------------
$quote = [char]0x201C + "All experiments ran inside a GPU-enabled VS Code Dev Container (Docker Desktop 28.1, CUDA 12.5 image) to ensure full reproducibility." + [char]0x201D

$anchor = $d.Content
$found = $anchor.Find.Execute("Advanced Audit Policy was set", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos = $anchor.Start
$insertPoint = $d.Range($startPos, $startPos)
$insertPoint.InsertBefore($quote)

$newSentenceRange = $d.Range($startPos, $startPos + $quote.Length)
$newSentenceRange.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# Change 2: collapse the three runs that spell out " /set /", "category:*",
# " /" (with a proofErr gramStart/gramEnd pair wrapping "category:*") into a
# single run reading " /set /category:* /". A plain Find/Replace across the
# whole span merges the identically-formatted runs and drops the now
# enclosed proofErr markers automatically.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(" /set /category:* /", $true, $false, $false, $false, $false, $true, 1, $false, " /set /category:* /", 2)

# ---------------------------------------------------------------------------
# Change 3a: merge "Supplementary " + "ZIP" (separated by a proofErr
# gramStart marker) into a single run "Supplementary ZIP", while leaving the
# following "  " (two-space) run untouched and separate.
#
# A direct Find/Replace over "Supplementary ZIP" would also merge in the
# next "  " run because it has identical run formatting and ends up
# adjacent after the edit. To avoid that, we briefly toggle off Bold on the
# "  " run so it no longer matches its neighbour's formatting, perform the
# merge, then restore Bold on it.
# ---------------------------------------------------------------------------
$locate = $d.Content
$null = $locate.Find.Execute("ZIP  Include", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$spacerStart = $locate.Start + 3
$spacerEnd = $spacerStart + 2
$spacer = $d.Range($spacerStart, $spacerEnd)
$spacer.Bold = 0

$null = $d.Content.Find.Execute("Supplementary ZIP", $true, $false, $false, $false, $false, $true, 1, $false, "Supplementary ZIP", 2)

$spacerAgain = $d.Range($spacerStart, $spacerEnd)
$spacerAgain.Bold = 1

# ---------------------------------------------------------------------------
# Change 3b: merge "Include" + " " (separated by a proofErr gramEnd marker)
# into a single run "Include ". The run immediately following (the Courier
# New "auditpol_backup.csv" run) has different formatting, so no
# unintended coalescing happens here.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Include ", $true, $false, $false, $false, $false, $true, 1, $false, "Include ", 2)

Write-Output "edit complete"
